$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.183.15'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '2.446.79'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.16'
$ws.Range("E5").Value = '  +2.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.06'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  +0.43%  '
$ws.Range("D9").Value = '2.441.64'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("E10").Value = '  +1.59%  '
$ws.Range("E11").Value = '  +2.76%  '
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.344'
$ws.Range("E13").Value = '  -2.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.52'
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("D16").Value = '2.868.35'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '62.035.58'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '2.432.93'
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.80'
$ws.Range("E19").Value = '  -2.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.21'
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.91'
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.11'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("E24").Value = '  -5.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.75'
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.10'
$ws.Range("E26").Value = '  +1.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '603.61'
$ws.Range("E27").Value = '  -2.91%  '
$ws.Range("D28").Value = '0.0₃0969'
$ws.Range("E28").Value = '  +0.67%  '
$ws.Range("E29").Value = '  +0.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.00'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("E32").Value = '  -1.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.91'
$ws.Range("E33").Value = '  +1.93%  '
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.89'
$ws.Range("E35").Value = '  -2.33%  '
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("E37").Value = '  -1.28%  '
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '153.87'
$ws.Range("E39").Value = '  +5.12%  '
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.29'
$ws.Range("E41").Value = '  +1.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '43.16'
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("E43").Value = '  -0.76%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '142.40'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.64'
$ws.Range("E47").Value = '  -1.49%  '
$ws.Range("E48").Value = '  +20.62%  '
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0521'
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.82'
$ws.Range("E51").Value = '  +0.08%  '
